$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# --- Slide-number placeholder -> "GRUPO 4 - ..." credits bar -------------
$ph = $s.Shapes.Item("Espaço Reservado para Número de Slide 6")

# Reposition / resize so the bar spans (almost) the full slide width.
$ph.Left = 7.091
$ph.Top = 498.7271
$ph.Width = 945.7405
$ph.Height = 28.75

$tf = $ph.TextFrame
$tr = $tf.TextRange
$tr.Text = "GRUPO 4" + " " + [char]0x2013 + " DANIEL SENA" + " " + [char]0x2013 + " GUSTAVO OLIVEIRA" + " " + [char]0x2013 + " LEANDRO BONETO" + " " + [char]0x2013 + " MIKKI DOS ANJOS" + " " + [char]0x2013 + " RODRIGO OLIVARES" + " " + [char]0x2013 + " SHELLY NADUDVARI"
$tr.Font.Size = 12
$tr.ParagraphFormat.Alignment = 2

$bold = $tr.Characters(1, 7)
$bold.Font.Bold = $true

# --- Remove the leftover "back" icon picture -------------------------------
$pic = $s.Shapes.Item("Picture 4")
$pic.Delete()
